$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of effort tracking to append after existing row 24 (previously last
# data row was row 25, which moves to row 27).

# Row 25: 03/07/2013, 1h, Revision manual
$ws.Range("A25").Value = 41458
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = "Revision manual"

# Row 26: 06/07/2013, 1h, Revision manual
$ws.Range("A26").Value = 41461
$ws.Range("B26").Value = 1
$ws.Range("D26").Value = "Revision manual"

# Row 27: 07/07/2013, 2.5h, SVN branch: gcc versus g++... (was row 25 previously)
$ws.Range("A27").Value = 41462
$ws.Range("B27").Value = 2.5
$ws.Range("D27").Value = "SVN branch: gcc versus g++. Revision of Makefile, support of Linux and Windows, modularization"

# Row 28: 08/07/2013, 2h, Revision manual
$ws.Range("A28").Value = 41463
$ws.Range("B28").Value = 2
$ws.Range("D28").Value = "Revision manual"

# Row 29: 09/07/2013, 1.5h, Implementation tc14
$ws.Range("A29").Value = 41464
$ws.Range("B29").Value = 1.5
$ws.Range("D29").Value = "Implementation tc14"

$ws.Range("A29").Select()
